$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row before row 41 (shifts old rows 41-1000 down to 42-1001) ---
$ws.Rows("41:41").Insert()

# Populate the newly inserted row 41 with blank "entry" line matching the
# style/layout of row 40 (CLICK FOR MERIT MAC 23), but with an empty label.
$ws.Range("A40:G40").Copy()
$ws.Range("A41:G41").PasteSpecial(-4122)
$ws.Range("E40").Copy()
$ws.Range("E41").PasteSpecial(-4122)

$ws.Range("B41").Value = ""
$ws.Range("C41").Value = ""
$ws.Range("D41").Value = 0
$ws.Range("E41").Value = 0
$ws.Range("F41").Value = ""

# Give D40 the same "readingOrder" style used by the C column (s=36) per the diff
$ws.Range("C40").Copy()
$ws.Range("D40").PasteSpecial(-4122)
$ws.Range("D40").Value = 0

# --- Row 42 (previously row 41, the "JUMLAH" subtotal line) ---
# Label changes from the "CLOSING BALANCE" text to the new "JUMLAH" text,
# and its formula gains a new D41-E41 term for the row we just inserted.
$ws.Range("B42").Value = "JUMLAH"
$ws.Range("F42").Formula = "=D34-E34+D35-E35+D36-E36+D37-E37+D38-E38+D39-E39+D41-E41+D40-E40"

# --- Row 43 (previously row 42, the "CLOSING BALANCE" total line) ---
$ws.Range("B43").Value = "CLOSING BALANCE FOR YEAR 2022 / BAKI PENUTUP TAHUN 2022"
$ws.Range("F43").Formula = "=F20+F32+F26+F42"

# --- Row 44 now blank (previously held the "End Of 2022 Statement" text) ---
$ws.Range("B44").Value = ""

# --- Row 45 keeps "End Of 2022 Statement" text ---
$ws.Range("B45").Value = "End Of 2022 Statement / Penyata 2022 Tamat "

# --- Row 46 gets the "May You Continue..." text (previously on row 45) ---
$ws.Range("B46").Value = "May You Continue to Rise Higher in 2023 "
$ws.Range("B45").Copy()
$ws.Range("B46").PasteSpecial(-4122)
$ws.Range("B46").Value = "May You Continue to Rise Higher in 2023 "

Write-Output "done"
